# Daily Report update: "Work Activities Vertical view completed"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PROJECT NAME
$ws.Range("B2").Value = "ELAC Campus Wide Duct Cleaning"

# DATE (kept as text, not an Excel date serial)
$ws.Range("F2").Value = "'08/13/2021"

# CONTRACT NO. (kept as text, not a number)
$ws.Range("B3").Value = "'4500289944"

# PROJECT ID
$ws.Range("F3").Value = 6300

# START TIME
$ws.Range("F4").Value = "09:32"

# WEATHER
$ws.Range("B5").Value = "Snowy"

# END TIME
$ws.Range("F5").Value = "09:49"

# Clear the two blank manpower/equipment rows (rows 10 and 11) entirely
$ws.Range("A10:F11").ClearContents()
